$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the header string in D1 from "Ten_year_flood_VIC_m3_sec" to "flood_discharge_VIC_m3_sec"
$ws.Range("D1").Value = "flood_discharge_VIC_m3_sec"

# Update selection to match the target state
$ws.Range("K1").Select()
